$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Price" (column D) and "Volume(1h)" (column E) cells hold plain text
# (not numbers) in the source data -- some of the new values look like
# numbers (e.g. "15.00", "0.05810") and would otherwise be silently
# coerced to a Double by the Value setter, losing significant trailing
# zeros. Force the cell to Text format first so the string is kept verbatim.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20.640.12"
$ws.Range("E2").Value = "  +2.76%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.475.50"
$ws.Range("E3").Value = "  +3.28%  "
$ws.Range("E4").Value = "  +0.59%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9884"
$ws.Range("E5").Value = "  -1.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "281.02"
$ws.Range("E6").Value = "  +2.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3732"
$ws.Range("E7").Value = "  +0.73%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3221"
$ws.Range("E8").Value = "  +4.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "41.87"
$ws.Range("E9").Value = "  +4.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.072"
$ws.Range("E10").Value = "  +6.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06758"
$ws.Range("E11").Value = "  +2.81%  "
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.701"
$ws.Range("E13").Value = "  +5.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.62"
$ws.Range("E14").Value = "  +8.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.338"
$ws.Range("E15").Value = "  +2.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.478.83"
$ws.Range("E16").Value = "  +3.02%  "
$ws.Range("E17").Value = "  +3.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.05810"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.58"
$ws.Range("E19").Value = "  -2.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9884"
$ws.Range("E20").Value = "  -1.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.749"
$ws.Range("E21").Value = "  +1.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.00"
$ws.Range("E22").Value = "  +3.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.28"
$ws.Range("E23").Value = "  +1.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.311"
$ws.Range("E24").Value = "  -0.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "20.760.86"
$ws.Range("E25").Value = "  +3.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.354"
$ws.Range("E26").Value = "  +2.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "138.80"
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.79"
$ws.Range("E28").Value = "  +5.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.642.24"
$ws.Range("E29").Value = "  +2.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "114.63"
$ws.Range("E30").Value = "  +4.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.995"
$ws.Range("E31").Value = "  +3.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.424"
$ws.Range("E32").Value = "  -0.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8549"
$ws.Range("E33").Value = "  -7.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.629"
$ws.Range("E34").Value = "  +24.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07896"
$ws.Range("E35").Value = "  +1.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06139"
$ws.Range("E36").Value = "  +8.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.996"
$ws.Range("E37").Value = "  +4.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.86"
$ws.Range("E38").Value = "  -5.21%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.864"
$ws.Range("E39").Value = "  -6.63%  "
$ws.Range("B40").Value = "Frax"
$ws.Range("C40").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9907"
$ws.Range("E40").Value = "  -1.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02092"
$ws.Range("E41").Value = "  +3.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.141"
$ws.Range("E42").Value = "  +1.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1921"
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5477"
$ws.Range("E44").Value = "  +2.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.64"
$ws.Range("E45").Value = "  +3.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.600"
$ws.Range("E46").Value = "  +1.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.18"
$ws.Range("E47").Value = "  +10.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5423"
$ws.Range("E48").Value = "  +5.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.842"
$ws.Range("E49").Value = "  +3.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.066"
$ws.Range("E50").Value = "  +1.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06481"
$ws.Range("E51").Value = "  +4.52%  "
